$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.501.72'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.45%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.504.11'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.05%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.45%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.59%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.495.63'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.611'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.47%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.199'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.26%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.645'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.65%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.70'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.87%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000305'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.42%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.42'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.064.07'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.21'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.31%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.470.27'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.49%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.502.78'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.21'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.60%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.119'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.35%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '539.60'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +11.44%  '

$ws.Range("E22").Value = '  -4.66%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '18.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -8.15%  '

$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.56'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.03%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.86'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '95.34'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.79%  '

$ws.Range("B27").Value = 'ImmutableX'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.96'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.13%  '

$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.07'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.37%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.06'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.60%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.95'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.68%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.31'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.43%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.54'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.16%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '64.49'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.15%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.113'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.01%  '

$ws.Range("B35").Value = 'Bittensor'
$ws.Range("C35").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '544.41'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.32%  '

$ws.Range("B36").Value = 'TheGraph'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.408'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.72%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.05'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.49%  '

$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.18%  '

$ws.Range("B39").Value = 'InjectiveProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.60'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.45%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0764'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.51%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.15'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.97%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.37'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.23%  '

$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.327.08'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.66%  '

$ws.Range("B44").Value = 'Kaspa'
$ws.Range("C44").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.132'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.96%  '

$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.51'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.49%  '

$ws.Range("B46").Value = 'ThetaToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.96'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.81%  '

$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0438'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.07%  '

$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.99'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.59%  '

$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.134'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.17%  '

$ws.Range("E50").Value = '  -0.09%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '137.01'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.94%  '
